# Generate Report for Handback
#
# This script mirrors a "localization handback" run: the zh-cn and de-de
# rows move from "Ready for handoff" to "Handed back: in sync with en-US",
# each language's "Latest Target File" / "Latest Handback File" columns get
# filled in with a hyperlink + the generated xliff name, the zh-cn /
# de-de "Latest Handback DateTime" timestamps are stamped, and a couple of
# report columns are widened so the new long file names are readable.

$wb = $excel.ActiveWorkbook

$hrefA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a92e8a194ea4d986cba6b9ead9572cecd26361f8/e2e/a.md"
$hrefB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a92e8a194ea4d986cba6b9ead9572cecd26361f8/e2e/b.md"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: status columns for both languages + widen the two
# language status columns (E, F) now that the text is longer.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("I2").Value = "a.md"
$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-21 14:45:19"

$wsZh.Range("I3").Value = "a.md"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-21 14:45:19"

$wsZh.Range("I2").Style = "Hyperlink"
$wsZh.Range("I3").Style = "Hyperlink"

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# Rebuild the hyperlinks collection in row/column order so link rIds line
# up the way Excel emits them (A2, I2, A3, I3).
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $hrefA, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $hrefA, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $hrefB, "", "", "b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $hrefA, "", "", "a.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("I2").Value = "a.md"
$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-21 14:45:26"

$wsDe.Range("I3").Value = "a.md"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-21 14:45:26"

$wsDe.Range("I2").Style = "Hyperlink"
$wsDe.Range("I3").Style = "Hyperlink"

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $hrefA, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $hrefA, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $hrefB, "", "", "b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $hrefA, "", "", "a.md")

Write-Host "Handback report generated"
